$wb = $excel.ActiveWorkbook

# Cell-level updates pulled from the scheduled market-price refresh.
# Each entry is either a plain value write, or a clear (when the
# upstream source no longer has a value for that cell).
$changes = @(
    @{ Sheet="ALC"; Cell="H80"; Value=6260.8 }
    @{ Sheet="ALC"; Cell="J80"; Value=10645.454 }
    @{ Sheet="ALC"; Cell="L80"; Value=31936.362 }
    @{ Sheet="ALC"; Cell="N80"; Value=-33932.362 }
    @{ Sheet="ALC"; Cell="H83"; Value=6260.8 }
    @{ Sheet="ALC"; Cell="J83"; Value=10645.454 }
    @{ Sheet="ALC"; Cell="L83"; Value=95809.086 }
    @{ Sheet="ALC"; Cell="N83"; Value=-105793.086 }
    @{ Sheet="ALC"; Cell="H88"; Value=2610 }
    @{ Sheet="ALC"; Cell="J88"; Value=2837.5 }
    @{ Sheet="ALC"; Cell="L88"; Value=2837.5 }
    @{ Sheet="ALC"; Cell="N88"; Value=-3649.5 }
    @{ Sheet="ALC"; Cell="H91"; Value=2610 }
    @{ Sheet="ALC"; Cell="J91"; Value=2837.5 }
    @{ Sheet="ALC"; Cell="L91"; Value=2837.5 }
    @{ Sheet="ALC"; Cell="N91"; Value=-5645.5 }
    @{ Sheet="ALC"; Cell="H98"; Value=3507.325 }
    @{ Sheet="ALC"; Cell="I98"; Value=2528.3547 }
    @{ Sheet="ALC"; Cell="J98"; Value=6879.3335 }
    @{ Sheet="ALC"; Cell="K98"; Value=2528.3547 }
    @{ Sheet="ALC"; Cell="L98"; Value=6879.3335 }
    @{ Sheet="ALC"; Cell="M98"; Value=-1030.3547 }
    @{ Sheet="ALC"; Cell="N98"; Value=-9875.3335 }
    @{ Sheet="ALC"; Cell="H113"; Value=3146.3333 }
    @{ Sheet="ALC"; Cell="I113"; Value=2487.5 }
    @{ Sheet="ALC"; Cell="J113"; Value=3585.5557 }
    @{ Sheet="ALC"; Cell="K113"; Value=2487.5 }
    @{ Sheet="ALC"; Cell="L113"; Value=3585.5557 }
    @{ Sheet="ALC"; Cell="M113"; Value=766.5 }
    @{ Sheet="ALC"; Cell="N113"; Value=-10093.5557 }
    @{ Sheet="ALC"; Cell="H122"; Value=3507.325 }
    @{ Sheet="ALC"; Cell="I122"; Value=2528.3547 }
    @{ Sheet="ALC"; Cell="J122"; Value=6879.3335 }
    @{ Sheet="ALC"; Cell="K122"; Value=7585.0641 }
    @{ Sheet="ALC"; Cell="L122"; Value=20638.0005 }
    @{ Sheet="ALC"; Cell="M122"; Value=-5135.0641 }
    @{ Sheet="ALC"; Cell="N122"; Value=-25538.0005 }
    @{ Sheet="ALC"; Cell="H132"; Value=1677.9688 }
    @{ Sheet="ALC"; Cell="I132"; Value=1286.8462 }
    @{ Sheet="ALC"; Cell="J132"; Value=3372.8333 }
    @{ Sheet="ALC"; Cell="K132"; Value=3860.5386 }
    @{ Sheet="ALC"; Cell="L132"; Value=10118.4999 }
    @{ Sheet="ALC"; Cell="M132"; Value=-1330.5386 }
    @{ Sheet="ALC"; Cell="N132"; Value=-15178.4999 }
    @{ Sheet="ALC"; Cell="H135"; Value=789.7895 }
    @{ Sheet="ALC"; Cell="I135"; Value=713.35486 }
    @{ Sheet="ALC"; Cell="J135"; Value=1128.2858 }
    @{ Sheet="ALC"; Cell="K135"; Value=6420.193740000001 }
    @{ Sheet="ALC"; Cell="L135"; Value=10154.5722 }
    @{ Sheet="ALC"; Cell="M135"; Value=-3885.193740000001 }
    @{ Sheet="ALC"; Cell="N135"; Value=-15224.5722 }
    @{ Sheet="ALC"; Cell="H138"; Value=2419.4426 }
    @{ Sheet="ALC"; Cell="I138"; Value=1359.878 }
    @{ Sheet="ALC"; Cell="J138"; Value=4591.55 }
    @{ Sheet="ALC"; Cell="K138"; Value=4079.634 }
    @{ Sheet="ALC"; Cell="L138"; Value=13774.65 }
    @{ Sheet="ALC"; Cell="M138"; Value=1060.366 }
    @{ Sheet="ALC"; Cell="N138"; Value=-24054.65 }
    @{ Sheet="ALC"; Cell="H141"; Value=7436.6665 }
    @{ Sheet="ALC"; Cell="I141"; Value=2703.8235 }
    @{ Sheet="ALC"; Cell="J141"; Value=27551.25 }
    @{ Sheet="ALC"; Cell="K141"; Value=8111.470499999999 }
    @{ Sheet="ALC"; Cell="L141"; Value=82653.75 }
    @{ Sheet="ALC"; Cell="M141"; Value=-2931.470499999999 }
    @{ Sheet="ALC"; Cell="N141"; Value=-93013.75 }
    @{ Sheet="ARM"; Cell="H32"; Value=14508.163 }
    @{ Sheet="ARM"; Cell="I32"; Value=14713.342 }
    @{ Sheet="ARM"; Cell="J32"; Value=12948.8 }
    @{ Sheet="ARM"; Cell="K32"; Value=14713.342 }
    @{ Sheet="ARM"; Cell="L32"; Value=12948.8 }
    @{ Sheet="ARM"; Cell="M32"; Value=-14426.342 }
    @{ Sheet="ARM"; Cell="N32"; Value=-13522.8 }
    @{ Sheet="ARM"; Cell="H43"; Value=333344260 }
    @{ Sheet="ARM"; Cell="J43"; Value=333344260 }
    @{ Sheet="ARM"; Cell="L43"; Value=333344260 }
    @{ Sheet="ARM"; Cell="N43"; Value=-333344886 }
    @{ Sheet="ARM"; Cell="H61"; Value=3852.7896 }
    @{ Sheet="ARM"; Cell="I61"; Value=3990.4443 }
    @{ Sheet="ARM"; Cell="J61"; Value=3728.9 }
    @{ Sheet="ARM"; Cell="K61"; Value=3990.4443 }
    @{ Sheet="ARM"; Cell="L61"; Value=3728.9 }
    @{ Sheet="ARM"; Cell="M61"; Value=-3778.4443 }
    @{ Sheet="ARM"; Cell="N61"; Value=-4152.9 }
    @{ Sheet="ARM"; Cell="H136"; Value=3852.7896 }
    @{ Sheet="ARM"; Cell="I136"; Value=3990.4443 }
    @{ Sheet="ARM"; Cell="J136"; Value=3728.9 }
    @{ Sheet="ARM"; Cell="K136"; Value=11971.3329 }
    @{ Sheet="ARM"; Cell="L136"; Value=11186.7 }
    @{ Sheet="ARM"; Cell="M136"; Value=-9421.332900000001 }
    @{ Sheet="ARM"; Cell="N136"; Value=-16286.7 }
    @{ Sheet="BSM"; Cell="H47"; Value=0 }
    @{ Sheet="BSM"; Cell="J47"; Value=0 }
    @{ Sheet="BSM"; Cell="L47"; Value=0 }
    @{ Sheet="BSM"; Cell="N47"; Clear=$true }
    @{ Sheet="BSM"; Cell="H59"; Value=79975 }
    @{ Sheet="BSM"; Cell="J59"; Value=79975 }
    @{ Sheet="BSM"; Cell="L59"; Value=79975 }
    @{ Sheet="BSM"; Cell="N59"; Value=-81669 }
    @{ Sheet="BSM"; Cell="H86"; Value=127300.94 }
    @{ Sheet="BSM"; Cell="I86"; Value=2937.625 }
    @{ Sheet="BSM"; Cell="J86"; Value=251664.25 }
    @{ Sheet="BSM"; Cell="K86"; Value=2937.625 }
    @{ Sheet="BSM"; Cell="L86"; Value=251664.25 }
    @{ Sheet="BSM"; Cell="M86"; Value=-1814.625 }
    @{ Sheet="BSM"; Cell="N86"; Value=-253910.25 }
    @{ Sheet="BSM"; Cell="H89"; Value=127300.94 }
    @{ Sheet="BSM"; Cell="I89"; Value=2937.625 }
    @{ Sheet="BSM"; Cell="J89"; Value=251664.25 }
    @{ Sheet="BSM"; Cell="K89"; Value=14688.125 }
    @{ Sheet="BSM"; Cell="L89"; Value=1258321.25 }
    @{ Sheet="BSM"; Cell="M89"; Value=-9072.125 }
    @{ Sheet="BSM"; Cell="N89"; Value=-1269553.25 }
    @{ Sheet="BSM"; Cell="H104"; Value=29710 }
    @{ Sheet="BSM"; Cell="J104"; Value=29710 }
    @{ Sheet="BSM"; Cell="L104"; Value=29710 }
    @{ Sheet="BSM"; Cell="N104"; Value=-36698 }
    @{ Sheet="CRP"; Cell="H50"; Value=35630 }
    @{ Sheet="CRP"; Cell="J50"; Value=35630 }
    @{ Sheet="CRP"; Cell="L50"; Value=35630 }
    @{ Sheet="CRP"; Cell="N50"; Value=-36880 }
    @{ Sheet="CRP"; Cell="H53"; Value=37728 }
    @{ Sheet="CRP"; Cell="J53"; Value=37728 }
    @{ Sheet="CRP"; Cell="L53"; Value=37728 }
    @{ Sheet="CRP"; Cell="N53"; Value=-38942 }
    @{ Sheet="CRP"; Cell="H132"; Value=367212.94 }
    @{ Sheet="CRP"; Cell="I132"; Value=501739.75 }
    @{ Sheet="CRP"; Cell="J132"; Value=3990.6 }
    @{ Sheet="CRP"; Cell="K132"; Value=1505219.25 }
    @{ Sheet="CRP"; Cell="L132"; Value=11971.8 }
    @{ Sheet="CRP"; Cell="M132"; Value=-1502689.25 }
    @{ Sheet="CRP"; Cell="N132"; Value=-17031.8 }
    @{ Sheet="CUL"; Cell="H101"; Value=9523.728 }
    @{ Sheet="CUL"; Cell="J101"; Value=9523.728 }
    @{ Sheet="CUL"; Cell="L101"; Value=28571.184 }
    @{ Sheet="CUL"; Cell="N101"; Value=-33439.18399999999 }
    @{ Sheet="CUL"; Cell="H131"; Value=2833.3867 }
    @{ Sheet="CUL"; Cell="I131"; Value=11500 }
    @{ Sheet="CUL"; Cell="J131"; Value=1651.5758 }
    @{ Sheet="CUL"; Cell="K131"; Value=34500 }
    @{ Sheet="CUL"; Cell="L131"; Value=4954.7274 }
    @{ Sheet="CUL"; Cell="M131"; Value=-29460 }
    @{ Sheet="CUL"; Cell="N131"; Value=-15034.7274 }
    @{ Sheet="GSM"; Cell="H117"; Value=40290 }
    @{ Sheet="GSM"; Cell="J117"; Value=40290 }
    @{ Sheet="GSM"; Cell="L117"; Value=40290 }
    @{ Sheet="GSM"; Cell="N117"; Value=-47174 }
    @{ Sheet="LTW"; Cell="H22"; Value=1025.125 }
    @{ Sheet="LTW"; Cell="I22"; Value=900 }
    @{ Sheet="LTW"; Cell="J22"; Value=1066.8334 }
    @{ Sheet="LTW"; Cell="K22"; Value=900 }
    @{ Sheet="LTW"; Cell="L22"; Value=1066.8334 }
    @{ Sheet="LTW"; Cell="M22"; Value=-605 }
    @{ Sheet="LTW"; Cell="N22"; Value=-1656.8334 }
    @{ Sheet="LTW"; Cell="H27"; Value=1025.125 }
    @{ Sheet="LTW"; Cell="I27"; Value=900 }
    @{ Sheet="LTW"; Cell="J27"; Value=1066.8334 }
    @{ Sheet="LTW"; Cell="K27"; Value=900 }
    @{ Sheet="LTW"; Cell="L27"; Value=1066.8334 }
    @{ Sheet="LTW"; Cell="M27"; Value=-793 }
    @{ Sheet="LTW"; Cell="N27"; Value=-1280.8334 }
    @{ Sheet="LTW"; Cell="H50"; Value=50000 }
    @{ Sheet="LTW"; Cell="J50"; Value=0 }
    @{ Sheet="LTW"; Cell="L50"; Value=0 }
    @{ Sheet="LTW"; Cell="N50"; Clear=$true }
    @{ Sheet="LTW"; Cell="H130"; Value=24421.5 }
    @{ Sheet="LTW"; Cell="J130"; Value=24421.5 }
    @{ Sheet="LTW"; Cell="L130"; Value=24421.5 }
    @{ Sheet="LTW"; Cell="N130"; Value=-34461.5 }
    @{ Sheet="LTW"; Cell="H136"; Value=16834952 }
    @{ Sheet="LTW"; Cell="I136"; Value=22728498 }
    @{ Sheet="LTW"; Cell="J136"; Value=627702.3 }
    @{ Sheet="LTW"; Cell="K136"; Value=68185494 }
    @{ Sheet="LTW"; Cell="L136"; Value=1883106.9 }
    @{ Sheet="LTW"; Cell="M136"; Value=-68182944 }
    @{ Sheet="LTW"; Cell="N136"; Value=-1888206.9 }
    @{ Sheet="WVR"; Cell="H45"; Value=0 }
    @{ Sheet="WVR"; Cell="J45"; Value=0 }
    @{ Sheet="WVR"; Cell="L45"; Value=0 }
    @{ Sheet="WVR"; Cell="N45"; Clear=$true }
    @{ Sheet="WVR"; Cell="H98"; Value=0 }
    @{ Sheet="WVR"; Cell="J98"; Value=0 }
    @{ Sheet="WVR"; Cell="L98"; Value=0 }
    @{ Sheet="WVR"; Cell="N98"; Clear=$true }
    @{ Sheet="WVR"; Cell="H136"; Value=2170.2104 }
    @{ Sheet="WVR"; Cell="I136"; Value=2218.3333 }
    @{ Sheet="WVR"; Cell="J136"; Value=1989.75 }
    @{ Sheet="WVR"; Cell="K136"; Value=6654.999899999999 }
    @{ Sheet="WVR"; Cell="L136"; Value=5969.25 }
    @{ Sheet="WVR"; Cell="M136"; Value=-4104.999899999999 }
    @{ Sheet="WVR"; Cell="N136"; Value=-11069.25 }
)

$wsCache = @{}

foreach ($change in $changes) {
    $sheetName = $change.Sheet
    if (-not $wsCache.ContainsKey($sheetName)) {
        $wsCache[$sheetName] = $wb.Worksheets.Item($sheetName)
    }
    $ws = $wsCache[$sheetName]
    $rng = $ws.Range($change.Cell)
    if ($change.ContainsKey("Clear") -and $change.Clear) {
        $rng.ClearContents()
    } else {
        $rng.Value = $change.Value
    }
}

Write-Host "Applied $($changes.Count) cell updates across $($wsCache.Count) sheets."
